$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. First paragraph: set alignment to left (adds <w:jc w:val="left"/>)
# ------------------------------------------------------------------
$d.Paragraphs(1).Range.ParagraphFormat.Alignment = 0   # wdAlignParagraphLeft

# ------------------------------------------------------------------
# 2. Drop the old "_GoBack" bookmark - it is going to be recreated
#    right after the simplified CSS sentence below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3. Simplify the CSS requirement sentence: drop the part about a
#    separate CSS file per page ("... e outro arquivo CSS
#    específico ... historia.html, historia.css.") and keep only the
#    "single estilo.css" requirement, then put the "_GoBack" bookmark
#    right after the new sentence (mirrors where Word leaves it after
#    the last edit).
# ------------------------------------------------------------------
$range = $d.Content
$range.Find.Execute(
    "Deve conter um arquivo CSS geral de nome estilo.css e outro arquivo CSS específico para cada página que você desenvolver. Ex. historia.html, historia.css.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deve conter um arquivo CSS geral de nome estilo.css.Z",
    2
) | Out-Null

# "Z" above is a one-character placeholder inserted right where the
# bookmark must end up. Shrink the range down to just that character...
$range.MoveStart(1, 52)

# ...anchor the bookmark on it (a non-empty range), then erase the
# placeholder through the bookmark's own Range object so the bookmark
# collapses to a zero-length mark in the correct spot.
$d.Bookmarks.Add("_GoBack", $range)
$d.Bookmarks("_GoBack").Range.Text = ""
